$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 26234.822
$ws.Range("J17").Value = 26234.822
$ws.Range("L17").Value = 78704.466
$ws.Range("N17").Value = -79040.466
$ws.Range("H40").Value = 39179.906
$ws.Range("I40").Value = 73581.82000000001
$ws.Range("J40").Value = 1337.8
$ws.Range("K40").Value = 73581.82000000001
$ws.Range("L40").Value = 1337.8
$ws.Range("M40").Value = -73406.82000000001
$ws.Range("N40").Value = -1687.8
$ws.Range("H64").Value = 3734.6667
$ws.Range("I64").Value = 3734.6667
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3734.6667
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -3486.6667
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 3734.6667
$ws.Range("I67").Value = 3734.6667
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3734.6667
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -2876.6667
$ws.Range("N67").ClearContents()
$ws.Range("H100").Value = 2960.05
$ws.Range("I100").Value = 2130.8333
$ws.Range("J100").Value = 4203.875
$ws.Range("K100").Value = 2130.8333
$ws.Range("L100").Value = 4203.875
$ws.Range("M100").Value = -1589.8333
$ws.Range("N100").Value = -5285.875
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960
$ws.Range("H134").Value = 35390
$ws.Range("J134").Value = 35390
$ws.Range("L134").Value = 35390
$ws.Range("N134").Value = -45530
$ws.Range("H135").Value = 1229.0605
$ws.Range("I135").Value = 979.6957
$ws.Range("J135").Value = 1802.6
$ws.Range("K135").Value = 8817.2613
$ws.Range("L135").Value = 16223.4
$ws.Range("M135").Value = -6282.2613
$ws.Range("N135").Value = -21293.4
$ws.Range("H137").Value = 2877.0715
$ws.Range("I137").Value = 3081
$ws.Range("J137").Value = 2129.3333
$ws.Range("K137").Value = 9243
$ws.Range("L137").Value = 6387.999899999999
$ws.Range("M137").Value = -6693
$ws.Range("N137").Value = -11487.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 725.62
$ws.Range("I32").Value = 720.0204
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 720.0204
$ws.Range("L32").Value = 1000
$ws.Range("M32").Value = -433.0204
$ws.Range("N32").Value = -1574

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1095.5
$ws.Range("I94").Value = 966.4286
$ws.Range("K94").Value = 966.4286
$ws.Range("M94").Value = -515.4286
$ws.Range("H105").Value = 1490.3572
$ws.Range("I105").Value = 1370
$ws.Range("J105").Value = 1744.4445
$ws.Range("K105").Value = 1370
$ws.Range("L105").Value = 1744.4445
$ws.Range("M105").Value = 377
$ws.Range("N105").Value = -5238.4445
$ws.Range("H130").Value = 29503.4
$ws.Range("J130").Value = 29503.4
$ws.Range("L130").Value = 29503.4
$ws.Range("N130").Value = -39543.4
$ws.Range("H132").Value = 15389.5
$ws.Range("J132").Value = 15389.5
$ws.Range("L132").Value = 15389.5
$ws.Range("N132").Value = -25509.5
$ws.Range("H134").Value = 1772.2759
$ws.Range("I134").Value = 1107.84
$ws.Range("K134").Value = 3323.52
$ws.Range("M134").Value = -788.5199999999995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1874.8667
$ws.Range("I16").Value = 1401.6666
$ws.Range("K16").Value = 1401.6666
$ws.Range("M16").Value = -1114.6666
$ws.Range("H41").Value = 4178.625
$ws.Range("I41").Value = 875
$ws.Range("J41").Value = 7482.25
$ws.Range("K41").Value = 875
$ws.Range("L41").Value = 7482.25
$ws.Range("M41").Value = -447
$ws.Range("N41").Value = -8338.25
$ws.Range("H62").Value = 6602.5
$ws.Range("I62").Value = 4205
$ws.Range("J62").Value = 9000
$ws.Range("K62").Value = 4205
$ws.Range("L62").Value = 9000
$ws.Range("M62").Value = -3581
$ws.Range("N62").Value = -10248
$ws.Range("H64").Value = 32333.334
$ws.Range("J64").Value = 32333.334
$ws.Range("L64").Value = 32333.334
$ws.Range("N64").Value = -32829.334
$ws.Range("H65").Value = 6602.5
$ws.Range("I65").Value = 4205
$ws.Range("J65").Value = 9000
$ws.Range("K65").Value = 21025
$ws.Range("L65").Value = 45000
$ws.Range("M65").Value = -17905
$ws.Range("N65").Value = -51240
$ws.Range("H67").Value = 32333.334
$ws.Range("J67").Value = 32333.334
$ws.Range("L67").Value = 32333.334
$ws.Range("N67").Value = -34049.334
$ws.Range("H68").Value = 39831.668
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 39831.668
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 39831.668
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -41329.668
$ws.Range("H69").Value = 17500
$ws.Range("I69").Value = 5000
$ws.Range("K69").Value = 5000
$ws.Range("M69").Value = -4251
$ws.Range("H71").Value = 39831.668
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 39831.668
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 119495.004
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -126983.004
$ws.Range("H72").Value = 17500
$ws.Range("I72").Value = 5000
$ws.Range("K72").Value = 15000
$ws.Range("M72").Value = -11256
$ws.Range("H105").Value = 3562.3076
$ws.Range("I105").Value = 4262
$ws.Range("J105").Value = 3125
$ws.Range("K105").Value = 4262
$ws.Range("L105").Value = 3125
$ws.Range("M105").Value = -2515
$ws.Range("N105").Value = -6619
$ws.Range("H113").Value = 1874.8667
$ws.Range("I113").Value = 1401.6666
$ws.Range("K113").Value = 1401.6666
$ws.Range("M113").Value = 768.3334
$ws.Range("H132").Value = 1948
$ws.Range("I132").Value = 1480.5454
$ws.Range("K132").Value = 4441.6362
$ws.Range("M132").Value = -1911.6362
$ws.Range("H134").Value = 2369.1155
$ws.Range("I134").Value = 1239.85
$ws.Range("J134").Value = 6133.3335
$ws.Range("K134").Value = 3719.55
$ws.Range("L134").Value = 18400.0005
$ws.Range("M134").Value = -1184.55
$ws.Range("N134").Value = -23470.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -2002
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -10008
$ws.Range("N83").ClearContents()
$ws.Range("H132").Value = 3729.5625
$ws.Range("I132").Value = 3726
$ws.Range("J132").Value = 3736.3635
$ws.Range("K132").Value = 11178
$ws.Range("L132").Value = 11209.0905
$ws.Range("M132").Value = -8648
$ws.Range("N132").Value = -16269.0905
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("H132").Value = 2471.6177
$ws.Range("I132").Value = 1521.4
$ws.Range("J132").Value = 5111.1113
$ws.Range("K132").Value = 4564.200000000001
$ws.Range("L132").Value = 15333.3339
$ws.Range("M132").Value = -2034.200000000001
$ws.Range("N132").Value = -20393.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 11806.472
$ws.Range("I132").Value = 2352.1282
$ws.Range("J132").Value = 38143.57
$ws.Range("K132").Value = 7056.3846
$ws.Range("L132").Value = 114430.71
$ws.Range("M132").Value = -4526.3846
$ws.Range("N132").Value = -119490.71
